# Update the cryptos list: Price (D) and Volume(1h) (E) columns
# for rows 2-51, per the Thu Dec 28 15:14:58 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.758.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.44%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.371.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.74%  "

# Row 4
$ws.Range("E4").Value = "  -0.32%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.69%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.637"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.61%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("E9").Value = "  +0.10%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.25%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0922"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.08%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.93%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.77%  "

# Row 14
$ws.Range("E14").Value = "  +0.70%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.88%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.727.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.65%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.367.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.76%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.673.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.51%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.14%  "

# Row 20
$ws.Range("E20").Value = "  -1.29%  "

# Row 21
$ws.Range("E21").Value = "  +7.14%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.51%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "273.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.96%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.26%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.47%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.71%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.37%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.61%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.16%  "

# Row 31
$ws.Range("E31").Value = "  -1.85%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0903"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.15%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.01%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.58%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.132"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.12%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.00%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0358"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.72%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.92%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.84%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.105"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.14%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.52%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +51.05%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.228"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.98%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "68.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.02%  "

# Row 45
$ws.Range("E45").Value = "  +0.00%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "116.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.36%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.55%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.18%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.01%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.610.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.22%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.15%  "
